$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.159824728965759
$ws.Range("B1").Value = 2.406866312026978
$ws.Range("C1").Value = -1
$ws.Range("D1").Value = 2.386054039001465
$ws.Range("E1").Value = 1.230137825012207
